# update layout master barang pelanggan
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old placeholder barcode text ("1234KDFH") with a distinct
# barcode value per row.
$ws.Range("L2").Value = "BARCODE1"
$ws.Range("L3").Value = "BARCODE2"
$ws.Range("L4").Value = "BARCODE3"
$ws.Range("L5").Value = "BARCODE4"

# The price/barcode columns (I:L) on the data rows used a one-off style that
# duplicated the plain bordered style used by the rest of the table (e.g.
# A2:H5). Re-apply that shared format so they match the surrounding cells.
[void]$ws.Range("A2").Copy()
$ws.Range("I2:L5").PasteSpecial(-4122) | Out-Null

# Widen the BARCODE column now that it holds the longer BARCODEn text.
$ws.Columns.Item(12).ColumnWidth = 11.02

# Move the active selection to K7.
$ws.Range("K7").Select() | Out-Null
